# Simulated Wild Card round and logged it
# Updates stat totals on the "Rushing" and "Receiving" sheets to reflect
# the outcome of a simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2 - J.Jackson
$rushing.Range("C2").Value = 15
$rushing.Range("D2").Value = 9
$rushing.Range("F2").Value = 14

# Row 3 - J.Kelley
$rushing.Range("C3").Value = 125
$rushing.Range("D3").Value = 66
$rushing.Range("E3").Value = 14
$rushing.Range("F3").Value = 46

# Row 4 - L.Rountree
$rushing.Range("C4").Value = 32

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 - J.Jackson
$receiving.Range("C2").Value = 79
$receiving.Range("D2").Value = 60
$receiving.Range("E2").Value = 6
$receiving.Range("G2").Value = 15
$receiving.Range("H2").Value = 13

# Row 3 - L.Rountree
$receiving.Range("C3").Value = 13

# Row 6 - K.Allen
$receiving.Range("C6").Value = 125
$receiving.Range("D6").Value = 85

# Row 7 - G.Nabers
$receiving.Range("C7").Value = 80
$receiving.Range("D7").Value = 51
$receiving.Range("E7").Value = 35
$receiving.Range("F7").Value = 14
$receiving.Range("G7").Value = 20
$receiving.Range("H7").Value = 8

# Row 8 - J.Palmer
$receiving.Range("C8").Value = 31
$receiving.Range("D8").Value = 21
$receiving.Range("E8").Value = 4
$receiving.Range("F8").Value = 2

# Row 9 - J.Guyton
$receiving.Range("C9").Value = 35
$receiving.Range("D9").Value = 22
$receiving.Range("G9").Value = 7
$receiving.Range("H9").Value = 4

# Row 11 - J.Cook
$receiving.Range("C11").Value = 56
$receiving.Range("D11").Value = 34
$receiving.Range("E11").Value = 13
$receiving.Range("F11").Value = 6
$receiving.Range("G11").Value = 8

# Row 14 - T.McKitty
$receiving.Range("C14").Value = 7
$receiving.Range("G14").Value = 1
